$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace a unique whole-paragraph text with new text.
# ---------------------------------------------------------------------------
function Replace-Text($oldText, $newText) {
    $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Output "WARNING: text not found for replacement: $oldText"
    }
}

# ---------------------------------------------------------------------------
# Helper: given the text of an existing (anchor) paragraph, insert one or
# more new bullet paragraphs immediately after it, in order. Uses the
# paragraph's numeric Index (not Next()/stale object refs) so it also works
# when the anchor is the very last paragraph in the document.
# ---------------------------------------------------------------------------
function Insert-ParasAfter($anchorText, $newTexts) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText)
    if (-not $found) {
        Write-Output "WARNING: anchor not found for insertion: $anchorText"
        return
    }
    $idx = $rng.Paragraphs(1).Index
    foreach ($t in $newTexts) {
        $d.Paragraphs($idx).Range.InsertParagraphAfter()
        $idx = $idx + 1
        $d.Paragraphs($idx).Range.Text = $t
    }
}

# ===========================================================================
# 1. CORE COMPETENCIES -- collapse three detailed paragraphs into one summary
# ===========================================================================
$d.Paragraphs(8).Range.Delete()
$d.Paragraphs(7).Range.Delete()
Replace-Text "Research and Analytics: Survey Methodology: Design, sampling, weighting, longitudinal analysis • Statistical Analysis: Regression modeling, clustering, segmentation, machine learning • Geospatial Analysis: Spatial clustering, boundary estimation, demographic mapping • Data Visualization: Tableau, PowerBI, d3.js, Matplotlib, Seaborn, choropleth mapping • Research Management: Team leadership, methodology design, stakeholder communication" "Research and Analytics • Programming and Development • Data Infrastructure"

# ===========================================================================
# 2. RESEARCH DIRECTOR - Progressive Change Campaign Committee
# ===========================================================================
Replace-Text "• Managed critical research operations for political campaigns" "• Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls"
Replace-Text "• Conducted comprehensive polling and demographic analysis" "• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren"
Replace-Text "• Developed strategic recommendations based on data analysis" "• Built tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver"
Replace-Text "• Led research team in support of progressive political initiatives" "• Designed survey deployment system facilitating thousands of simultaneous phone surveys"
Insert-ParasAfter "• Designed survey deployment system facilitating thousands of simultaneous phone surveys" @(
    "• Significantly increased data collection efficiency through automated calling infrastructure",
    "• Managed comprehensive research operations for progressive political initiatives and candidates"
)

# ===========================================================================
# 3. SOFTWARE ENGINEER - Salsa Labs, Inc.
# ===========================================================================
Replace-Text "• Developed software solutions for political campaigns and advocacy groups" "• Maintained and extended entire geospatial analysis and reporting tools for Java-based CRM system"
Replace-Text "• Built web applications for voter engagement and campaign management" "• Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers"
Replace-Text "• Integrated third-party APIs and data sources for campaign tools" "• Built geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill"
Replace-Text "• Collaborated with political strategists to translate requirements into technical solutions" "• Integrated mapping and visualization tools for political campaign data analysis"
Insert-ParasAfter "• Integrated mapping and visualization tools for political campaign data analysis" @(
    "• Collaborated with political strategists to translate geospatial requirements into technical solutions"
)

# ===========================================================================
# 4. INTERIM TECHNOLOGY MANAGER - The Praxis Project
# ===========================================================================
Replace-Text "• Integrated technology solutions within organizational frameworks for social justice organizations" "• Assisted in search for full-time CTO while performing all programmatic technology roles for multi-million dollar organization"
Replace-Text "• Developed data management systems for community organizing efforts" "• Made all technology decisions and practices for massive multinational non-governmental organization"
Replace-Text "• Provided technical training and support to nonprofit staff" "• Wrote comprehensive frameworks for internal and external technology audits"
Replace-Text "• Built custom applications for community engagement and advocacy" "• Trained beneficiaries on spatial and Census data analysis for public health research"
Insert-ParasAfter "• Trained beneficiaries on spatial and Census data analysis for public health research" @(
    "• Trained NGO staff in web development using Drupal, PHP, and MySQL",
    "• Managed technology infrastructure supporting community health initiatives across multiple countries"
)

# ===========================================================================
# 5. PROGRAMMER - Lake Research Partners
# ===========================================================================
Replace-Text "• Developed data analysis tools for political polling and research" "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party"
Replace-Text "• Built statistical models for voter behavior analysis" "• Developed system that later became the Polling Consortium Database at The Analyst Institute"
Replace-Text "• Created data visualization tools for research presentations" "• Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections"
Replace-Text "• Supported senior researchers with technical analysis and reporting" "• Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle"
Insert-ParasAfter "• Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle" @(
    "• Pioneered integration of advanced mapping techniques into standard reports including choropleths and hexagonal grid maps",
    "• Developed innovative approaches to visualizing demographic and market data for enhanced client understanding"
)

# ===========================================================================
# 6. FIELD DIRECTOR - The Feldman Group
# ===========================================================================
Replace-Text "• Managed field operations for political campaigns and research projects" "• Administered all quantitative and qualitative research operations ensuring reporting accuracy"
Replace-Text "• Developed data collection and management systems for field work" "• Managed comprehensive survey fielding for multi-million dollar research firm"
Replace-Text "• Trained field staff on data collection protocols and quality control" "• Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings"
Replace-Text "• Analyzed field data to inform campaign strategy and research findings" "• Created custom reports and data visualizations based on specific client requirements"
Insert-ParasAfter "• Created custom reports and data visualizations based on specific client requirements" @(
    "• Introduced mapping and geospatial analysis into standard reporting procedures",
    "• Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL"
)

# ===========================================================================
# 7. New TECHNICAL SKILLS section at the very end of the document
# ===========================================================================
Insert-ParasAfter "• Redistricting analysis used in court cases with rigorous methodology and expert testimony" @(
    "TECHNICAL SKILLS",
    "RESEARCH AND ANALYTICS Survey Methodology; Statistical Analysis; Geospatial Analysis; Data Visualization; Research Management",
    "PROGRAMMING AND DEVELOPMENT Python; JVM Languages; Web Technologies; Database Languages; Statistical Computing",
    "DATA INFRASTRUCTURE Cloud Platforms; Big Data; Databases; Geospatial; DevOps"
)
$rng = $d.Content
$rng.Find.Execute("TECHNICAL SKILLS") | Out-Null
$rng.Paragraphs(1).Style = "Heading 2"

Write-Output "Edit complete"
